# -----------------------------------------------------------------------
# "Extra Wishlist" sheet reshuffle:
#   - Rows re-sorted/grouped (no longer a plain price sort of all items);
#     blank "GAP" separator rows inserted between groups (col E = "GAP").
#   - The last group of 4 rows additionally gets col E = "Y".
#   - All 21 original product rows are preserved verbatim, just moved.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old data block (values AND formatting) - rows are being
# reordered, and several rows need to change from "Hyperlink"-styled link
# cells to plain style (and vice versa), so a clean slate is simplest.
$ws.Range("A2:E22").Clear()

# Each entry: RowNumber, Name, Image, Link, LinkStyle("Hyperlink"/"Normal"), Price, ReservedMarker
$data = @(
    @(2, $null, $null, $null, "Normal", $null, "GAP")
    ,@(3, "Apple Watch", "https://store.storeimages.cdn-apple.com/1/as-images.apple.com/is/MDRY4ref_VW_34FR+watch-case-42-aluminum-rosegold-nc-s10_VW_34FR+watch-face-42-aluminum-rosegold-s10_VW_34FR?wid=5120&hei=3280&bgc=fafafa&trim=1&fmt=p-jpg&qlt=80&.v=aWR3eURSMEloeW8rRlJjeisrTUVuSjNqbENGcEFFTlJVaXJwL2VzdGxEM3hoaFEySTlrZSsrbEhqczIxajV3MzBnRnJNVU50TjZtbUN3ZjYvVksrL1FtTHNhRmRKQkE2OXJHRVEvRTFWdXRGd1ZOdUM2RkNJOUpIeklzU1JHbVNWNEM0ZHZzMzMzajJiOGZuamRiTE9rUVN3R3VxZWhYYXgwOHljYmZFMXBocmMyRTN3NCt6QkoxaUdRb0FBay9VYktGTHdENW9lYUFnak5pcy9ReEdDV1huc2RkMVJhUksxbWsxK1ZQMHVvWVdkR0c1QU9HSFhRTXZyMjRWUXMzUA", "https://www.apple.com/ch-de/shop/buy-watch/apple-watch/42mm-gps-ros%C3%A9gold-aluminium-mandarine-sportarmband-sm", "Normal", 400, $null)
    ,@(4, "Apple iPad Air 11 2025 (M3)", "https://static01.galaxus.com/productimages/1/9/1/3/1/6/9/0/2/7/6/0/5/5/9/9/7/1/56e6de25-8e03-4f19-b997-d12d17f2a578_cropped.jpg_720.avif", "https://www.digitec.ch/en/s1/product/apple-ipad-air-11-2025-m3-wlan-only-11-256-gb-space-grey-tablets-55620779", "Hyperlink", 700, $null)
    ,@(5, $null, $null, $null, "Normal", $null, "GAP")
    ,@(6, "Sigma 24-70mm f/2.8 DG DN II Art", "https://static01.galaxus.com/productimages/3/4/1/4/6/9/1/6/1/2/1/6/6/2/2/1/8/4/9/682182dc-69c0-4397-b129-6fef5b7c1b49_cropped.jpg_2880.avif", "https://www.digitec.ch/en/s1/product/sigma-24-70mm-f28-dg-dn-ii-art-sony-e-full-size-lenses-45881601", "Normal", 1160, $null)
    ,@(7, "Voigtländer Macro Apo-Lanthar 110 mm f/2.5", "https://static01.galaxus.com/productimages/2/9/0/4/5/5/5/9/voigtlander.jpg_480.avif", "https://www.digitec.ch/en/s1/product/voigtlaender-macro-apo-lanthar-110-mm-f25-sony-e-full-size-lenses-11921001", "Hyperlink", 1200, $null)
    ,@(8, "CITY OF DARKNESS – SPECIAL REPRINT", "https://pictures.abebooks.com/isbn/9781873200896-us.jpg", "https://www.abebooks.com/City-Darkness-Limited-Edition-Life-Kowloon/31292330335/bd", "Normal", 160, $null)
    ,@(9, "Sony FE 40mm f/2.5 G", "https://static01.galaxus.com/productimages/1/9/9/3/2/3/5/8/4/3/6/8/0/4/6/0/0/8/2/01981e19-c81d-77c1-a3c2-f578d3601927_2880.avif", "https://www.digitec.ch/en/s1/product/sony-fe-40mm-f25-g-sony-e-full-size-lenses-15382655", "Normal", 520, $null)
    ,@(10, $null, $null, $null, "Normal", $null, "GAP")
    ,@(11, "Sony Alpha 7RV", "https://static01.galaxus.com/productimages/1/9/8/0/7/4/3/7/2/2/3/8/6/1/4/8/5/4/8/935761d6-8268-40f6-83dd-a55b34216305_cropped.jpg_720.avif", "https://www.digitec.ch/en/s1/product/sony-alpha-7rv-61-mpx-full-frame-cameras-22823231", "Hyperlink", 3800, $null)
    ,@(12, "Sony FE 35mm f/1.4 GM", "https://static01.galaxus.com/productimages/4/1/4/3/9/0/3/8/1%20Digitec_A-Mid%203.jpg_2880.avif", "https://www.digitec.ch/en/s1/product/sony-fe-35mm-f14-gm-sony-e-full-size-lenses-14596684", "Normal", 1236, $null)
    ,@(13, "Sigma 100-400mm f/5.0-6.3 DG DN OS, Sony E", "https://static01.galaxus.com/productimages/3/5/9/9/6/4/4/5/PPhoto_100_400_5_6.3_dgdn_c020_Lmt_horizontal_tripodsocket.jpg_720.avif", "https://www.digitec.ch/en/s1/product/sigma-100-400mm-f50-63-dg-dn-os-sony-e-sony-e-full-size-lenses-13366507", "Hyperlink", 870, $null)
    ,@(14, $null, $null, $null, "Normal", $null, "GAP")
    ,@(15, "Deuter Forest Fox 10", "https://static01.galaxus.com/productimages/3/3/3/3/4/0/2/0/0/1/7/4/8/9/2/2/0/5/2/58cd5acc-8462-46b1-8934-2a09baf72ca7_cropped.jpg_720.avif", "https://www.galaxus.ch/en/s8/product/deuter-forest-fox-10-10-l-backpacks-42902192", "Hyperlink", 45, $null)
    ,@(16, "Deuter Forest Fox 10", "https://static01.galaxus.com/productimages/1/6/1/0/9/8/0/5/7/2/6/1/7/7/2/6/8/4/5/118d9d2b-617e-4c6c-ba5b-b3e0fac55598_cropped.jpg_720.avif", "https://www.galaxus.ch/en/s8/product/deuter-waldfuchs-10-10-l-backpacks-42902191", "Hyperlink", 45, $null)
    ,@(17, $null, $null, $null, "Normal", $null, "GAP")
    ,@(18, "Synology DS923+", "https://static01.galaxus.com/productimages/3/6/3/4/3/0/7/0/5/5/9/8/3/8/0/5/9/0/9/1a5b9310-192f-4565-8e87-b6a9163262f9_cropped.jpg_2880.avif", "https://www.digitec.ch/en/s1/product/synology-ds923-0-tb-nas-23108534", "Hyperlink", 540, $null)
    ,@(19, "SVS SB-1000 Pro", "https://www.mueller-spring.ch/media/image/product/388254/lg/lautsprecher-subwoofer-svs-sb-1000-pro-black-ash~2.webp", "https://www.mueller-spring.ch/lautsprecher/subwoofer/svs-sb-1000-pro-black-ash", "Hyperlink", 750, $null)
    ,@(20, "SSD", "https://static01.galaxus.com/productimages/2/8/4/6/9/5/7/5/9/5/8/7/7/8/8/8/6/3/9/c9f8d1b6-3490-4d83-adc8-b60c5dea8143_cropped.jpg_720.avif", "https://www.digitec.ch/en/s1/product/wd-black-sn850x-4000-gb-m2-2280-ssd-21635315", "Hyperlink", 250, $null)
    ,@(21, "LEGO 10297 Boutique Hotel", "https://www.lego.com/cdn/cs/set/assets/bltf5232acff5ba0c26/10297_alt1.png?format=webply&fit=bounds&quality=75&width=800&height=800&dpr=1", "https://www.lego.com/en-ch/product/boutique-hotel-10297", "Normal", 240, $null)
    ,@(22, "LEGO 10350 Tudor Corner", "https://www.lego.com/cdn/cs/set/assets/bltc38b5911d77d1733/10350_Prod.png?format=webply&fit=bounds&quality=75&width=1200&height=1200&dpr=1", "https://www.lego.com/en-ch/product/tudor-corner-10350", "Normal", 240, $null)
    ,@(23, "Aeron Chair", "https://images.hermanmiller.group/m/49e8a45a115fae9a/W-HM_AER_61524.png?trim=auto&trim-sd=1&blend-mode=darken&blend=fafafa&bg=fafafa&auto=format&w=1000&q=70&h=1000", "https://store.hermanmiller.com/office-chairs-aeron/aeron-chair/100249179.html?lang=en_US&sku=100249179", "Hyperlink", 1700, $null)
    ,@(24, $null, $null, $null, "Normal", $null, "GAP")
    ,@(25, "Fujifilm X100VI", "https://static01.galaxus.com/productimages/3/4/2/1/3/9/9/5/7/5/9/5/8/9/1/2/3/5/1/f55d3e63-3122-4c75-b336-1095e95d55f0_cropped.jpg_720.avif", "https://www.digitec.ch/en/s1/product/fujifilm-x100vi-23-mm-4020-mpx-aps-c-dx-cameras-42944935", "Hyperlink", 1600, "Y")
    ,@(26, "Sony FE 90mm f/2.8 Makro G OSS", "https://static01.galaxus.com/productimages/8/5/5/8/2/8/2/4/1/9/2/2/6/6/9/0/5/2/9/9d13d7b3-49a5-41fb-a402-2a1a04b162c9_cropped.jpg_720.avif", "https://www.digitec.ch/en/s1/product/sony-fe-90mm-f28-makro-g-oss-sony-e-full-size-lenses-3528209", "Hyperlink", 750, "Y")
    ,@(27, "The Orient Express Train", "https://www.lego.com/cdn/cs/set/assets/blt0c893d9aba078a6a/21344_boxprod_v39.png?format=webply&fit=bounds&quality=75&width=800&height=800&dpr=1", "https://www.lego.com/en-ch/product/the-orient-express-train-21344?CMP=AFC-AffiliateUK-JF2DY0n6uhs-3146997-1539294-10", "Hyperlink", 330, "Y")
    ,@(28, "Dire Straits Live 1978 - 1992 LP", "https://superdeluxeedition.com/wp-content/uploads/2023/09/vinyl_spread-1.jpg", "https://www.amazon.de/dp/B0C92W29LD?tag=sdepcwde-21&linkCode=ogi&th=1&psc=1", "Normal", 220, "Y")
)

foreach ($row in $data) {
    $r      = $row[0]
    $name   = $row[1]
    $image  = $row[2]
    $link   = $row[3]
    $style  = $row[4]
    $price  = $row[5]
    $marker = $row[6]

    if ($name -ne $null)   { $ws.Cells.Item($r, 1).Value = $name }
    if ($image -ne $null)  { $ws.Cells.Item($r, 2).Value = $image }
    if ($link -ne $null) {
        $ws.Cells.Item($r, 3).Value = $link
        $ws.Cells.Item($r, 3).Style = $style
    }
    if ($price -ne $null)  { $ws.Cells.Item($r, 4).Value = $price }
    if ($marker -ne $null) { $ws.Cells.Item($r, 5).Value = $marker }
}

# Row 3 (first data row) was selected (whole row) when the file was saved.
$ws.Rows(3).Select()
